$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet (this
#    carries over all formatting / column widths / header styling) and place
#    it right after "总计", i.e. before the old "2022-Q3" sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# --- Row 2: brand-new fund entry for this quarter ---------------------------
# (leading "'" keeps these numeric-looking values stored as text, matching
# the source data convention used throughout this workbook)
$q4Sheet.Range("A2").Value2 = 0
$q4Sheet.Range("B2").Value2 = "'161810"
$q4Sheet.Range("C2").Value2 = "银华内需精选混合（LOF）"
$q4Sheet.Range("D2").Value2 = "'23.02"
$q4Sheet.Range("E2").Value2 = "'93.80"
$q4Sheet.Range("F2").Value2 = "'4.69"
$q4Sheet.Range("G2").Value2 = "'1.0796"
$q4Sheet.Range("H2").Value2 = 9

# --- Row 3: the fund already tracked in 2022-Q3, with this quarter's figures
$q4Sheet.Range("A3").Value2 = 1
$q4Sheet.Range("B3").Value2 = "'161838"
$q4Sheet.Range("C3").Value2 = "银华创业板两年定期开放混合"
$q4Sheet.Range("D3").Value2 = "'4.71"
$q4Sheet.Range("E3").Value2 = "'95.23"
$q4Sheet.Range("F3").Value2 = "'6.89"
$q4Sheet.Range("G3").Value2 = "'0.3245"
$q4Sheet.Range("H3").Value2 = 7

# New row needs the same look (bold / bordered / centered index cell) as the
# rest of column A; copy that formatting down from row 2.
$q4Sheet.Range("A2").Copy()
$q4Sheet.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row for 2022-Q4 above the
#    existing quarters, pushing everything else down, then renumber the
#    sequential index column (A).
# ---------------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert()

$totalSheet.Range("A2").Value2 = 0
$totalSheet.Range("B2").Value2 = "2022-Q4"
$totalSheet.Range("C2").Value2 = 2
$totalSheet.Range("D2").Value2 = 1.4

# The freshly inserted row inherits the header row's (bold/centered) look;
# reset it to the plain data-row formatting used by every other data row,
# then restore the distinct index-column (A) styling.
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Re-sequence the index column for the remaining (shifted) rows.
$totalSheet.Range("A3").Value2 = 1
$totalSheet.Range("A4").Value2 = 2

# Append the new final row for 2022-Q1 (previously the last quarter shown).
$totalSheet.Range("A5").Value2 = 3
$totalSheet.Range("B5").Value2 = "2022-Q1"
$totalSheet.Range("C5").Value2 = 1
$totalSheet.Range("D5").Value2 = 0.6

# Match formatting of the new row to the row above it.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Restore the originally-active tab ("总计"), since inserting / copying
#    sheets above shifts the active-tab focus as a side effect.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("总计").Activate()
